$wb = $excel.ActiveWorkbook

function Set-HeaderStyle($range) {
    $range.Font.Bold = $true
    $range.Borders.LineStyle = 1
    $range.HorizontalAlignment = -4108
    $range.VerticalAlignment = -4160
}

# ------------------------------------------------------------------
# 1. The existing "总计" sheet becomes the new "2022-Q1" quarterly
#    holdings sheet (same sheetId/position it already had), and a
#    brand-new "总计" summary sheet is appended right after it.
# ------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"
$total.Outline.SummaryRow = 1
$total.Outline.SummaryColumn = 1

# ------------------------------------------------------------------
# 2. Rebuild "2022-Q1" with the single fund-holding row supplied.
# ------------------------------------------------------------------
$q1.Cells.Clear()

Set-HeaderStyle($q1.Range("B1:H1"))
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

Set-HeaderStyle($q1.Range("A2"))
$q1.Range("A2").Value = 0
$q1.Range("B2").NumberFormat = "@"
$q1.Range("B2").Value = "001112"
$q1.Range("C2").Value = "东方红中国优势灵活配置混合"
$q1.Range("D2").NumberFormat = "@"
$q1.Range("D2").Value = "40.53"
$q1.Range("E2").NumberFormat = "@"
$q1.Range("E2").Value = "88.28"
$q1.Range("F2").NumberFormat = "@"
$q1.Range("F2").Value = "2.97"
$q1.Range("G2").NumberFormat = "@"
$q1.Range("G2").Value = "1.2037"
$q1.Range("H2").Value = 10

# ------------------------------------------------------------------
# 3. Populate the new "总计" sheet: same quarterly roll-up as before,
#    with a freshly inserted "2022-Q1" row on top and every other
#    row's running index shifted down by one.
# ------------------------------------------------------------------
Set-HeaderStyle($total.Range("B1:D1"))
$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$rows = @(
    @("2022-Q1", 1, 1.2),
    @("2021-Q4", 2, 3.21),
    @("2021-Q3", 2, 3.34),
    @("2021-Q2", 19, 16.89),
    @("2021-Q1", 13, 10.44),
    @("2020-Q4", 11, 11.4)
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $entry = $rows[$i]
    Set-HeaderStyle($total.Range("A$r"))
    $total.Range("A$r").Value = $i
    $total.Range("B$r").Value = $entry[0]
    $total.Range("C$r").Value = $entry[1]
    $total.Range("D$r").Value = $entry[2]
}
